{"js": "// Replace each arithmetic-problem cell's text with its updated value, in\n// reading order (row-major, 20 rows x 5 columns). The table's shape\n// (1 table, 20 rows, 5 cols, 100 cells, no merges) is unchanged by the\n// edit -- only the cell contents move/change -- so a straight positional\n// remap of every cell's text reproduces the diff exactly.\nconst newValues = [\n  \"68-25=\",\n  \"92-86=\",\n  \"29+47=\",\n  \"80-34=\",\n  \"72-35=\",\n  \"51+45=\",\n  \"14+2=\",\n  \"87+1=\",\n  \"39+6=\",\n  \"89-52=\",\n  \"2+1=\",\n  \"97-34=\",\n  \"59+23=\",\n  \"79-30=\",\n  \"88-60=\",\n  \"52-39=\",\n  \"60+28=\",\n  \"63-9=\",\n  \"62-35=\",\n  \"35+62=\",\n  \"31+49=\",\n  \"6+46=\",\n  \"66-60=\",\n  \"68-20=\",\n  \"49-42=\",\n  \"74-63=\",\n  \"84-0=\",\n  \"55+44=\",\n  \"62+5=\",\n  \"96-60=\",\n  \"31+10=\",\n  \"48-4=\",\n  \"30-21=\",\n  \"29+38=\",\n  \"82-52=\",\n  \"24-3=\",\n  \"65+14=\",\n  \"58-1=\",\n  \"89+8=\",\n  \"93-62=\",\n  \"91+7=\",\n  \"36+30=\",\n  \"66+33=\",\n  \"37+47=\",\n  \"92-47=\",\n  \"89-85=\",\n  \"82-13=\",\n  \"90-87=\",\n  \"7+37=\",\n  \"31-9=\",\n  \"64-36=\",\n  \"30+40=\",\n  \"36+1=\",\n  \"10+62=\",\n  \"35+19=\",\n  \"34+36=\",\n  \"71-67=\",\n  \"99-43=\",\n  \"24+46=\",\n  \"77+1=\",\n  \"81+8=\",\n  \"33+19=\",\n  \"6+28=\",\n  \"80+3=\",\n  \"96-89=\",\n  \"53+3=\",\n  \"81-65=\",\n  \"20+34=\",\n  \"38-34=\",\n  \"8+60=\",\n  \"73-71=\",\n  \"59-26=\",\n  \"45+44=\",\n  \"75-61=\",\n  \"80-60=\",\n  \"64-7=\",\n  \"20+35=\",\n  \"13+46=\",\n  \"19+60=\",\n  \"74-50=\",\n  \"70-14=\",\n  \"3+7=\",\n  \"92-50=\",\n  \"24+55=\",\n  \"4+51=\",\n  \"44+8=\",\n  \"51+30=\",\n  \"13-9=\",\n  \"85-8=\",\n  \"76-8=\",\n  \"12+46=\",\n  \"52-29=\",\n  \"30+26=\",\n  \"87-53=\",\n  \"93-92=\",\n  \"15+24=\",\n  \"8+28=\",\n  \"74+11=\",\n  \"80-3=\",\n  \"4+42=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = newValues.length / rowCount;\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    table.getCell(r, c).value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-problem cell's text with its updated value, in\n# reading order (row-major, 20 rows x 5 columns). The table's shape\n# (1 table, 20 rows, 5 cols, 100 cells, no merges) is unchanged by the\n# edit -- only the cell contents move/change -- so a straight positional\n# remap of every cell's text reproduces the diff exactly.\n$newValues = @(\n  \"68-25=\",\n  \"92-86=\",\n  \"29+47=\",\n  \"80-34=\",\n  \"72-35=\",\n  \"51+45=\",\n  \"14+2=\",\n  \"87+1=\",\n  \"39+6=\",\n  \"89-52=\",\n  \"2+1=\",\n  \"97-34=\",\n  \"59+23=\",\n  \"79-30=\",\n  \"88-60=\",\n  \"52-39=\",\n  \"60+28=\",\n  \"63-9=\",\n  \"62-35=\",\n  \"35+62=\",\n  \"31+49=\",\n  \"6+46=\",\n  \"66-60=\",\n  \"68-20=\",\n  \"49-42=\",\n  \"74-63=\",\n  \"84-0=\",\n  \"55+44=\",\n  \"62+5=\",\n  \"96-60=\",\n  \"31+10=\",\n  \"48-4=\",\n  \"30-21=\",\n  \"29+38=\",\n  \"82-52=\",\n  \"24-3=\",\n  \"65+14=\",\n  \"58-1=\",\n  \"89+8=\",\n  \"93-62=\",\n  \"91+7=\",\n  \"36+30=\",\n  \"66+33=\",\n  \"37+47=\",\n  \"92-47=\",\n  \"89-85=\",\n  \"82-13=\",\n  \"90-87=\",\n  \"7+37=\",\n  \"31-9=\",\n  \"64-36=\",\n  \"30+40=\",\n  \"36+1=\",\n  \"10+62=\",\n  \"35+19=\",\n  \"34+36=\",\n  \"71-67=\",\n  \"99-43=\",\n  \"24+46=\",\n  \"77+1=\",\n  \"81+8=\",\n  \"33+19=\",\n  \"6+28=\",\n  \"80+3=\",\n  \"96-89=\",\n  \"53+3=\",\n  \"81-65=\",\n  \"20+34=\",\n  \"38-34=\",\n  \"8+60=\",\n  \"73-71=\",\n  \"59-26=\",\n  \"45+44=\",\n  \"75-61=\",\n  \"80-60=\",\n  \"64-7=\",\n  \"20+35=\",\n  \"13+46=\",\n  \"19+60=\",\n  \"74-50=\",\n  \"70-14=\",\n  \"3+7=\",\n  \"92-50=\",\n  \"24+55=\",\n  \"4+51=\",\n  \"44+8=\",\n  \"51+30=\",\n  \"13-9=\",\n  \"85-8=\",\n  \"76-8=\",\n  \"12+46=\",\n  \"52-29=\",\n  \"30+26=\",\n  \"87-53=\",\n  \"93-92=\",\n  \"15+24=\",\n  \"8+28=\",\n  \"74+11=\",\n  \"80-3=\",\n  \"4+42=\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $tbl.Cell($r, $c).Range.Text = $newValues[$i]\n    $i++\n  }\n}\n"}
